$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.830.44"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "3.547.65"
$ws.Range("E3").Value = "  +1.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "617.09"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.04"
$ws.Range("E6").Value = "  +3.53%  "
$ws.Range("D7").Value = "3.544.71"
$ws.Range("E7").Value = "  +1.44%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.488"
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.145"
$ws.Range("E10").Value = "  +5.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.31"
$ws.Range("E11").Value = "  +5.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.437"
$ws.Range("E12").Value = "  +3.58%  "
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.94"
$ws.Range("E14").Value = "  +4.33%  "
$ws.Range("D15").Value = "4.146.40"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("D16").Value = "3.551.67"
$ws.Range("E16").Value = "  +1.51%  "
$ws.Range("D17").Value = "67.845.77"
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.117"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.76"
$ws.Range("E19").Value = "  +5.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.90"
$ws.Range("E20").Value = "  +5.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.93"
$ws.Range("E21").Value = "  +10.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "453.21"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.639"
$ws.Range("E23").Value = "  +2.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.83"
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000130"
$ws.Range("E25").Value = "  +1.18%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.48"
$ws.Range("E26").Value = "  +3.26%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "3.681.39"
$ws.Range("E27").Value = "  +1.26%  "
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.04"
$ws.Range("E29").Value = "  +8.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.56"
$ws.Range("E30").Value = "  +2.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.68"
$ws.Range("E31").Value = "  +5.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.169"
$ws.Range("E32").Value = "  +2.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.38"
$ws.Range("E34").Value = "  +4.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.99"
$ws.Range("E35").Value = "  +0.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.92"
$ws.Range("E36").Value = "  +3.51%  "
$ws.Range("D37").Value = "3.538.88"
$ws.Range("E37").Value = "  +1.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.23"
$ws.Range("E38").Value = "  +3.15%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.35"
$ws.Range("E39").Value = "  +6.41%  "
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "178.90"
$ws.Range("E41").Value = "  +4.40%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0915"
$ws.Range("E43").Value = "  +5.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.55"
$ws.Range("E44").Value = "  +2.56%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "30.76"
$ws.Range("E45").Value = "  +15.11%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.892"
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.33"
$ws.Range("E47").Value = "  +6.56%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.58"
$ws.Range("E48").Value = "  +0.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.63"
$ws.Range("E49").Value = "  +3.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.78"
$ws.Range("E50").Value = "  +3.28%  "
$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.260"
$ws.Range("E51").Value = "  +6.20%  "
